# The "escalas" header in C1 was renamed to the singular "escala",
# and the active selection moved from D2 back to C2 (with the view
# scrolled back so column A is visible again, i.e. no frozen/scrolled
# topLeftCell override).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "escala"

$ws.Activate()
$ws.Range("C2").Select()
